$d = $word.ActiveDocument
$d.Content.Find.Execute("has to", $true, $false, $false, $false, $false, $true, 1, $false, "must", 2)
